$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# G2: area of first segment
$ws.Range("G2").Formula = "=(D2-0)*B2/100"

# G3: single (non-shared in XML but same formula pattern) segment area
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

# G4:G15 shared formula for area segments
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# H2: total area
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# J2 / K2: summary values
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# Update selected cell to match target
$ws.Range("H8").Select()
